$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.110.80'
$ws.Range('E2').Value = '  -0.18%  '

$ws.Range('D3').Value = '1.653.84'
$ws.Range('E3').Value = '  -0.29%  '

$ws.Range('E4').Value = '  -0.24%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '217.51'
$ws.Range('E5').Value = '  +0.56%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5243'
$ws.Range('E6').Value = '  +0.71%  '

$ws.Range('E7').Value = '  -0.24%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2601'
$ws.Range('E8').Value = '  -1.22%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06337'
$ws.Range('E9').Value = '  +0.81%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '20.30'
$ws.Range('E10').Value = '  -2.34%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07792'
$ws.Range('E11').Value = '  +0.91%  '

$ws.Range('B12').Value = 'WrappedEther'
$ws.Range('C12').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D12').Value = '1.694.15'
$ws.Range('E12').Value = '  +2.19%  '

$ws.Range('B13').Value = 'Polkadot'
$ws.Range('C13').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.498'
$ws.Range('E13').Value = '  +1.57%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.5461'
$ws.Range('E14').Value = '  +0.74%  '

$ws.Range('D15').Value = '0.0₅8182'
$ws.Range('E15').Value = '  +0.46%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '65.31'
$ws.Range('E16').Value = '  +1.37%  '

$ws.Range('D17').Value = '26.110.08'
$ws.Range('E17').Value = '  -0.34%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '1.001'
$ws.Range('E18').Value = '  -0.34%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '4.571'
$ws.Range('E19').Value = '  -1.09%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '191.04'
$ws.Range('E20').Value = '  -0.39%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '10.05'
$ws.Range('E21').Value = '  -0.12%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.008'
$ws.Range('E22').Value = '  -0.86%  '

$ws.Range('E23').Value = '  -0.38%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '142.09'
$ws.Range('E24').Value = '  +1.73%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.1236'
$ws.Range('E25').Value = '  +0.62%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.229'
$ws.Range('E26').Value = '  +0.67%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '16.10'
$ws.Range('E27').Value = '  +0.50%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.429'
$ws.Range('E28').Value = '  +1.95%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.05891'
$ws.Range('E29').Value = '  -1.28%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.277'
$ws.Range('E30').Value = '  +0.65%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.505'
$ws.Range('E31').Value = '  -1.58%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.241'
$ws.Range('E32').Value = '  -0.52%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.582'
$ws.Range('E33').Value = '  -1.53%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.9482'
$ws.Range('E34').Value = '  -1.66%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.780'
$ws.Range('E35').Value = '  +0.27%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.409'
$ws.Range('E36').Value = '  -0.65%  '

$ws.Range('E37').Value = '  +0.06%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01617'
$ws.Range('E38').Value = '  +1.20%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '5.808'
$ws.Range('E39').Value = '  -3.12%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.8472'
$ws.Range('E40').Value = '  -1.15%  '

$ws.Range('E41').Value = '  -0.19%  '

$ws.Range('D42').Value = '1.025.96'
$ws.Range('E42').Value = '  +1.30%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '102.59'
$ws.Range('E43').Value = '  +2.14%  '

$ws.Range('D44').Value = '1.798.34'
$ws.Range('E44').Value = '  -0.11%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '57.06'
$ws.Range('E45').Value = '  +0.59%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.003'
$ws.Range('E46').Value = '  -0.35%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.4305'
$ws.Range('E47').Value = '  +2.61%  '

$ws.Range('E48').Value = '  -0.12%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '7.857'
$ws.Range('E49').Value = '  -1.42%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.468'
$ws.Range('E50').Value = '  +1.31%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.09682'
$ws.Range('E51').Value = '  -0.53%  '
